$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the redundant "Mani Bhavan Gandhi Sanghralaya" constraint row
# (row 7) - this shifts the following rows up and drops the now-unused
# shared string.
$ws.Rows.Item(7).Delete()

# Leave the selection on the row that took its place.
$ws.Range("A7").Select()
